$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = "telway"
$ws.Range("A5").Value = "northern TV and vacuum"

$ws.Range("A6").Select()
